# Reschedule Validations Update and Completion
#
# Adds a "Date" column (C) to the LoanAction sheet: a bold "Date" header
# in C2 (matching the existing Term/Interest header style) and a
# quote-prefixed, date-formatted text value "2/10/2014" in C3 (matching
# the TestCase Name/Value layout already used on the sheet), then
# autosizes the new column.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("LoanAction")

# Header cell C2: "Date", bold like the existing A2/B2 headers (Term/Interest)
$ws.Range("C2").Value = "Date"
$ws.Range("C2").Font.Bold = $true

# Data cell C3: literal text "2/10/2014" (leading apostrophe forces text,
# not an auto-converted date serial), formatted with a date display format
$ws.Range("C3").Value = "'2/10/2014"
$ws.Range("C3").NumberFormat = "mm-dd-yy"

# Resize column C to fit its new contents
$ws.Columns("C").AutoFit()
